# Add 2022-Q3 data:
#  1. Insert a new worksheet "2022-Q3" right after "总计" (before the
#     current "2022-Q2" sheet) and fill it with the quarterly holdings table.
#  2. Update the "总计" (summary) sheet so its top row now reports the
#     2022-Q3 numbers, cascading every older quarter down one row and
#     appending a new row for 2020-Q4 at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q3" sheet, inserted before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$col = 2
foreach ($h in $headers) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
    $col = $col + 1
}

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0,  "000311", "景顺长城沪深300指数增强A",           "46.75", "93.66", "2.17", "1.0145", 10),
    @(1,  "510160", "南方中证南方小康产业ETF",             "2.34",  "99.43", "7.75", "0.1814", 1),
    @(2,  "320011", "诺安中小盘精选混合",                   "3.02",  "73.05", "5.63", "0.1700", 2),
    @(3,  "160919", "大成产业升级股票（LOF）",              "3.07",  "84.12", "3.99", "0.1225", 8),
    @(4,  "159811", "博时中证5G产业50ETF",                  "2.04",  "98.58", "5.45", "0.1112", 5),
    @(5,  "320015", "诺安行业轮动混合",                     "1.11",  "73.94", "4.79", "0.0532", 4),
    @(6,  "009840", "西藏东财量化精选混合A",                "1.51",  "85.15", "3.34", "0.0504", 7),
    @(7,  "517960", "上投摩根中证沪港深科技100ETF",         "1.62",  "93.51", "2.70", "0.0437", 9),
    @(8,  "512220", "景顺长城中证科技传媒通信150ETF",       "2.35",  "98.48", "1.78", "0.0418", 10),
    @(9,  "006429", "诺安恒鑫混合",                         "0.64",  "70.94", "5.51", "0.0353", 4),
    @(10, "009841", "西藏东财量化精选混合C",                "0.67",  "85.15", "3.34", "0.0224", 7),
    @(11, "005326", "景顺长城泰恒回报灵活配置混合C",       "3.62",  "20.68", "0.57", "0.0206", 6),
    @(12, "003015", "中金沪深300指数增强A",                 "1.66",  "93.60", "1.21", "0.0201", 10),
    @(13, "517360", "华安中证沪港深科技100ETF",             "0.59",  "94.47", "2.66", "0.0157", 9),
    @(14, "517090", "国泰富时中国国企开放共赢ETF",         "0.60",  "91.75", "2.06", "0.0124", 9),
    @(15, "585001", "东吴中证新兴",                         "0.49",  "93.18", "1.90", "0.0093", 2),
    @(16, "510990", "工银瑞信中证180ESGETF",                "0.60",  "97.85", "1.49", "0.0089", 10),
    @(17, "010307", "西藏东财信息产业精选混合A",           "0.17",  "84.91", "4.87", "0.0083", 4),
    @(18, "005035", "银华信息科技量化优选股票A",           "0.19",  "89.32", "3.21", "0.0061", 3),
    @(19, "003579", "中金沪深300指数增强C",                 "0.47",  "93.60", "1.21", "0.0057", 10),
    @(20, "005325", "景顺长城泰恒回报灵活配置混合A",       "0.97",  "20.68", "0.57", "0.0055", 6),
    @(21, "014649", "永赢优质精选混合A",                     "0.09",  "80.97", "4.75", "0.0043", 4),
    @(22, "010308", "西藏东财信息产业精选混合C",           "0.08",  "84.91", "4.87", "0.0039", 4),
    @(23, "015679", "景顺长城沪深300指数增强C",             "0.06",  "93.66", "2.17", "0.0013", 10),
    @(24, "005036", "银华信息科技量化优选股票C",           "0.04",  "89.32", "3.21", "0.0013", 3),
    @(25, "014650", "永赢优质精选混合C",                     "0.00",  "80.97", "4.75", 0,        4)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
    $idxCell.Borders.Weight = 2

    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $sizeCell = $newSheet.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[3]

    $stockCell = $newSheet.Cells.Item($r, 5)
    $stockCell.NumberFormat = "@"
    $stockCell.Value = $row[4]

    $ratioCell = $newSheet.Cells.Item($r, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $row[5]

    # Last row's 持有市值 is a genuine number (0), every earlier row keeps it
    # as formatted text.
    $mvCell = $newSheet.Cells.Item($r, 7)
    if ($r -eq 27) {
        $mvCell.Value = $row[6]
    } else {
        $mvCell.NumberFormat = "@"
        $mvCell.Value = $row[6]
    }

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet with the new 2022-Q3 row, shifting
#    every other quarter down and appending 2020-Q4 at the bottom.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$summary = @(
    @(0, "2022-Q3", 26, 1.97),
    @(1, "2022-Q2", 12, 0.67),
    @(2, "2022-Q1", 16, 1.17),
    @(3, "2021-Q4", 13, 1),
    @(4, "2021-Q3", 20, 1.88),
    @(5, "2021-Q2", 25, 8.96),
    @(6, "2021-Q1", 13, 5.1),
    @(7, "2020-Q4", 6, 0.51)
)

$r = 2
foreach ($row in $summary) {
    $idxCell = $total.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
    $idxCell.Borders.Weight = 2

    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]

    $r = $r + 1
}
